$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-12-04"

# Update the header label for the current-year column (shared string / cell I1)
$ws.Range("I1").Value = "2022 (through 12-04)"

# Update December (row 13) and Total (row 14) values for the 2022 column (I)
$ws.Range("I13").Value = 16
$ws.Range("I14").Value = 1532
